# Apply "Atualizacoes 16 de janeiro de 2024" edits to the PROMIS-10 ValueSet workbook.

$wb = $excel.ActiveWorkbook

# 1. Rename the second worksheet tab.
$wsInclude = $wb.Worksheets.Item("Include from ")
$wsInclude.Name = "Include from PROMIS-10 Respon"

# 2. Update values on the "Metadata" sheet.
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.1 -> 0.0.0
$wsMeta.Range("B3").Value = "0.0.0"

# Title: ValueSet of PROMIS-10 Response Questionnaire -> PROMIS-10 Response Questionnaire
$wsMeta.Range("B5").Value = "PROMIS-10 Response Questionnaire"

# Experimental: (blank) -> "false" (must stay a text cell, not a Boolean,
# so write it through a formula that evaluates to text and paste the
# result in as a value - this avoids Excel's literal TRUE/FALSE
# autodetection that a direct .Value="false" assignment would trigger).
$scratch = $wsMeta.Range("D1")
$scratch.Formula = "=""false"""
$scratch.Copy()
$wsMeta.Range("B7").PasteSpecial(-4163)
$scratch.Clear()

# Date: 2023-11-21T19:08:35-03:00 -> 2024-01-11T13:00:00-03:00
$wsMeta.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# Description wording change
$wsMeta.Range("B12").Value = "ValueSet that defines the response values for the PROMIS-10 (Patient-Reported Outcomes Measurement Information System 10-item short form) Questionnaire."

# 3. Update the System URI value on the "Include from ..." sheet.
$wsInclude.Range("B34").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/PROMIS10VS"
